$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KY_PND_SEQ_TRANS value on row 2
$ws.Range("A2").Value = 63098612

# Update the comment text (shared string used by both row 2 and row 3)
$ws.Range("G2").Value = "QA automation test same supplier"

# Fill in previously-empty row 3 with a new record
$ws.Range("A3").Value = 61838590
$ws.Range("B3").Value = "Y"
$ws.Range("C3").Value = "N"
$ws.Range("D3").Value = "'07"
$ws.Range("E3").Value = "N"
$ws.Range("F3").Value = "WinkelJ"
$ws.Range("G3").Value = "QA automation test same supplier"

# Move the active selection to the newly added row
$ws.Range("A3").Select()
